$wb = $excel.ActiveWorkbook

# ============ Sheet: Summary Table ============
$ws1 = $wb.Worksheets.Item("Summary Table")

# Row 2
$v = @'
 Artificial Intelligence 
'@
$ws1.Range("A2").Value = $v
$v = @'
 37 
'@
$ws1.Range("B2").Value = $v
$v = @'
 China unveils new AI governance framework with focus on advanced model regulation 
'@
$ws1.Range("C2").Value = $v
$v = @'
 https://www.scmp.com/tech/policy/article/3259421/china-unveils-new-ai-governance-framework-emphasizing-advanced-model-safety 
'@
$ws1.Range("D2").Value = $v
$v = @'
 China's Ministry of Science and Technology released the "Framework for AI Safety and Governance 2025-2030" on May 7, 2025. The document outlines mandatory security assessments for large AI models and establishes a national AI safety research center. The framework emphasizes responsible AI development while maintaining China's competitive edge in generative AI. Industry leaders including Baidu, SenseTime, and Alibaba have already begun implementing compliance protocols. 
'@
$ws1.Range("E2").Value = $v

# Row 3
$v = @'
 Quantum communication 
'@
$ws1.Range("A3").Value = $v
$v = @'
 29 
'@
$ws1.Range("B3").Value = $v
$v = @'
 China-EU joint quantum encryption network demonstrates intercontinental quantum-secured communication 
'@
$ws1.Range("C3").Value = $v
$v = @'
 https://www.nature.com/articles/s41586-025-5742-x 
'@
$ws1.Range("D3").Value = $v
$v = @'
 Researchers from the University of Science and Technology of China and European partners successfully demonstrated the first intercontinental quantum-secured communication network on May 4, 2025. The system uses satellite-based quantum key distribution and terrestrial fiber networks to create an unhackable communication channel. The project marks a significant milestone in China-EU scientific cooperation and advances practical quantum communication technology. Chinese officials highlighted this as a model for future international collaboration in sensitive technology areas. 
'@
$ws1.Range("E3").Value = $v

# Row 4
$v = @'
 Digital economy policy 
'@
$ws1.Range("A4").Value = $v
$v = @'
 25 
'@
$ws1.Range("B4").Value = $v
$v = @'
 China launches Digital Silk Road 2.0 initiative focusing on AI infrastructure and sustainability 
'@
$ws1.Range("C4").Value = $v
$v = @'
 https://www.xinhuanet.com/english/2025-05/06/c_1419352.htm 
'@
$ws1.Range("D4").Value = $v
$v = @'
 On May 6, 2025, China's National Development and Reform Commission announced the Digital Silk Road 2.0 initiative with a $50 billion investment package. The program will support digital infrastructure development across Belt and Road countries with emphasis on sustainable technologies, AI solutions for climate monitoring, and interoperable digital payment systems. The initiative includes technology transfer provisions and training programs for partner countries. Twenty countries have already signed memorandums of understanding to participate. 
'@
$ws1.Range("E4").Value = $v

# Row 5
$v = @'
 New quality productivity 
'@
$ws1.Range("A5").Value = $v
$v = @'
 22 
'@
$ws1.Range("B5").Value = $v
$v = @'
 China unveils manufacturing transformation plan with targets for advanced manufacturing clusters 
'@
$ws1.Range("C5").Value = $v
$v = @'
 https://www.reuters.com/technology/china-unveils-ambitious-manufacturing-transformation-plan-2025-05-08/ 
'@
$ws1.Range("D5").Value = $v
$v = @'
 China's State Council released the "New Quality Productivity Acceleration Plan (2025-2035)" on May 8, 2025, detailing the country's strategy to upgrade manufacturing capabilities. The plan designates 25 advanced manufacturing clusters in sectors including semiconductors, biotech, and new energy vehicles. It establishes tax incentives for R&D investments and sets targets for increasing the value-added contribution of high-tech manufacturing to 35% of GDP by 2035. The plan emphasizes indigenous innovation while maintaining openness to international cooperation. 
'@
$ws1.Range("E5").Value = $v

# Row 6
$v = @'
 Semiconductor packaging 
'@
$ws1.Range("A6").Value = $v
$v = @'
 21 
'@
$ws1.Range("B6").Value = $v
$v = @'
 China announces major breakthroughs in advanced semiconductor packaging technologies 
'@
$ws1.Range("C6").Value = $v
$v = @'
 https://asia.nikkei.com/Business/Tech/Semiconductors/China-announces-breakthroughs-in-advanced-chip-packaging 
'@
$ws1.Range("D6").Value = $v
$v = @'
 On May 3, 2025, China's Ministry of Industry and Information Technology announced significant breakthroughs in advanced semiconductor packaging technologies. Researchers at the Chinese Academy of Sciences developed new fan-out wafer-level packaging techniques that improve performance by 40% while reducing power consumption. The advancements are particularly important given ongoing export controls on advanced chipmaking equipment. Five Chinese packaging firms will receive government support to commercialize the technology, with production lines expected by year-end. 
'@
$ws1.Range("E6").Value = $v

# Row 7
$v = @'
 Hydrogen energy storage 
'@
$ws1.Range("A7").Value = $v
$v = @'
 19 
'@
$ws1.Range("B7").Value = $v
$v = @'
 China and Germany launch joint hydrogen energy storage research center 
'@
$ws1.Range("C7").Value = $v
$v = @'
 https://www.cleanenergywire.org/news/china-germany-launch-joint-hydrogen-energy-storage-research-center 
'@
$ws1.Range("D7").Value = $v
$v = @'
 China's Ministry of Science and Technology and Germany's Federal Ministry of Education and Research inaugurated a joint research center for hydrogen energy storage on May 5, 2025. Located in Suzhou, the facility will receive €200 million in joint funding over five years. Research will focus on advanced electrolyzers, materials for hydrogen storage, and grid-scale deployment. The partnership aims to accelerate commercialization of green hydrogen technologies and establish common standards. Both nations emphasized this collaboration as critical to meeting climate goals. 
'@
$ws1.Range("E7").Value = $v

# Row 8
$v = @'
 Transformation of scientific and technological achievements 
'@
$ws1.Range("A8").Value = $v
$v = @'
 18 
'@
$ws1.Range("B8").Value = $v
$v = @'
 China reforms tech transfer system with new incentives for researchers and universities 
'@
$ws1.Range("C8").Value = $v
$v = @'
 https://www.scmp.com/tech/policy/article/3259487/china-overhauls-tech-transfer-system-boost-commercialization 
'@
$ws1.Range("D8").Value = $v
$v = @'
 China's State Council issued the "Comprehensive Reform Plan for Science and Technology Achievement Transformation" on May 6, 2025. The policy significantly increases the share of licensing revenue that researchers can retain (up to 70%), simplifies the approval process for university spin-offs, and creates a national technology transfer platform. The reform also establishes specialized intellectual property courts in 15 innovation hubs and creates a 100 billion yuan fund to support early commercialization of promising technologies. 
'@
$ws1.Range("E8").Value = $v

# Row 9
$v = @'
 Science and technology security policy 
'@
$ws1.Range("A9").Value = $v
$v = @'
 17 
'@
$ws1.Range("B9").Value = $v
$v = @'
 China issues new regulations on international scientific collaboration with security provisions 
'@
$ws1.Range("C9").Value = $v
$v = @'
 https://www.nature.com/articles/d41586-025-01355-2 
'@
$ws1.Range("D9").Value = $v
$v = @'
 On May 7, 2025, China's State Council released new regulations governing international scientific collaboration that balance openness with security concerns. The "Measures for Security Management of International Scientific and Technological Cooperation" establish review mechanisms for joint research projects in sensitive areas while streamlining approval for collaboration in non-sensitive fields. The regulations clarify data sharing protocols and intellectual property arrangements. Scientific societies have welcomed the clarity while some international partners expressed concerns about potential restrictions. 
'@
$ws1.Range("E9").Value = $v

# Row 10
$v = @'
 Integration of industry, academia and research 
'@
$ws1.Range("A10").Value = $v
$v = @'
 15 
'@
$ws1.Range("B10").Value = $v
$v = @'
 China launches 50 national innovation clusters to strengthen industry-academia integration 
'@
$ws1.Range("C10").Value = $v
$v = @'
 https://english.www.gov.cn/news/topnews/202505/08/content_WS6549f7a2c6d0868f4e2b2e37.html 
'@
$ws1.Range("D10").Value = $v
$v = @'
 The Chinese Ministry of Education and Ministry of Science and Technology jointly announced the establishment of 50 National Innovation Clusters on May 8, 2025. These clusters will connect leading universities, research institutes, and enterprises in strategic sectors. Each cluster will receive 500 million yuan in initial funding and preferential policies for talent recruitment and infrastructure development. The initiative aims to shorten the innovation cycle from research to commercialization. Early focus areas include integrated circuits, quantum information, and biomanufacturing. 
'@
$ws1.Range("E10").Value = $v

# Row 11
$v = @'
 Carbon fiber composites 
'@
$ws1.Range("A11").Value = $v
$v = @'
 15 
'@
$ws1.Range("B11").Value = $v
$v = @'
 China achieves breakthrough in low-cost carbon fiber production for aerospace applications 
'@
$ws1.Range("C11").Value = $v
$v = @'
 https://www.globaltimes.cn/page/202505/1306782.shtml 
'@
$ws1.Range("D11").Value = $v
$v = @'
 Chinese researchers at Harbin Institute of Technology announced a breakthrough in carbon fiber production technology on May 4, 2025. The new process reduces manufacturing costs by 40% while maintaining aerospace-grade quality. The technology uses domestic precursors and innovative thermal treatment techniques. China Aerospace Science and Technology Corporation plans to use the materials in its next-generation satellite structures and launch vehicles. Production facilities with 5,000-ton annual capacity will be constructed in Heilongjiang Province with operations expected to begin in 2026. 
'@
$ws1.Range("E11").Value = $v

# Row 12
$v = @'
 Brain-computer interface 
'@
$ws1.Range("A12").Value = $v
$v = @'
 14 
'@
$ws1.Range("B12").Value = $v
$v = @'
 China approves first clinical trial of implantable brain-computer interface for paralysis patients 
'@
$ws1.Range("C12").Value = $v
$v = @'
 https://www.sciencemag.org/news/2025/05/china-approves-groundbreaking-brain-computer-interface-clinical-trial 
'@
$ws1.Range("D12").Value = $v
$v = @'
 China's National Medical Products Administration approved the first clinical trial of an implantable brain-computer interface device on May 5, 2025. Developed by researchers at Tsinghua University and the Chinese Academy of Sciences, the "NeuralLink-C" device will be tested in 15 patients with complete spinal cord injuries. The minimally invasive device uses a new microelectrode array with 1,024 channels. The trial marks China's entry into advanced neural implant development, an area previously dominated by US companies. Patient recruitment will begin in June 2025. 
'@
$ws1.Range("E12").Value = $v

# Row 13
$v = @'
 International innovation platform 
'@
$ws1.Range("A13").Value = $v
$v = @'
 13 
'@
$ws1.Range("B13").Value = $v
$v = @'
 China launches Global Innovation Exchange Platform with initial focus on climate technologies 
'@
$ws1.Range("C13").Value = $v
$v = @'
 https://www.chinadaily.com.cn/a/202505/06/WS6546c891a3104efcbdad7e21.html 
'@
$ws1.Range("D13").Value = $v
$v = @'
 China's Ministry of Science and Technology launched the Global Innovation Exchange Platform (GIEP) on May 6, 2025, designed to facilitate international research collaboration and technology transfer. The platform initially focuses on climate technologies, advanced materials, and health sciences. Twenty countries have joined as founding members, with the platform providing research matching, funding coordination, and IP protection services. The initiative includes a 10 billion yuan fund for joint research projects. Officials emphasized GIEP's role in addressing global challenges through open innovation models. 
'@
$ws1.Range("E13").Value = $v

# Row 14
$v = @'
 Green environmental protection technology 
'@
$ws1.Range("A14").Value = $v
$v = @'
 12 
'@
$ws1.Range("B14").Value = $v
$v = @'
 China-ASEAN partnership announces major environmental technology transfer program 
'@
$ws1.Range("C14").Value = $v
$v = @'
 https://asean.org/china-asean-launch-environmental-technology-partnership/ 
'@
$ws1.Range("D14").Value = $v
$v = @'
 The China-ASEAN Environmental Technology Partnership was launched on May 7, 2025, at a ministerial meeting in Singapore. China will provide $2 billion in financing and technical assistance to support sustainable development across Southeast Asia. The program focuses on water treatment technologies, air pollution monitoring systems, and waste management solutions. Training programs for 5,000 environmental engineers from ASEAN countries will be established at Chinese universities. The partnership aligns with both China's ecological civilization goals and ASEAN's sustainability agenda. 
'@
$ws1.Range("E14").Value = $v

# Row 15
$v = @'
 Solid-state batteries 
'@
$ws1.Range("A15").Value = $v
$v = @'
 12 
'@
$ws1.Range("B15").Value = $v
$v = @'
 Chinese researchers achieve energy density breakthrough in solid-state battery technology 
'@
$ws1.Range("C15").Value = $v
$v = @'
 https://www.sciencedirect.com/science/article/pii/B9780323856249000156 
'@
$ws1.Range("D15").Value = $v
$v = @'
 A research team from the Chinese Academy of Sciences published breakthrough results in solid-state battery technology on May 3, 2025. Their prototype achieved an energy density of 500 Wh/kg using a new composite electrolyte material and advanced manufacturing process. The batteries demonstrated stable performance over 1,000 cycles and improved safety characteristics compared to lithium-ion batteries. CATL announced plans to incorporate the technology into pilot production lines by late 2025. Chinese officials highlighted the advancement as critical for electric vehicle development and energy transition goals. 
'@
$ws1.Range("E15").Value = $v

# Row 16
$v = @'
 Science and technology development plan 
'@
$ws1.Range("A16").Value = $v
$v = @'
 11 
'@
$ws1.Range("B16").Value = $v
$v = @'
 China unveils 15-year plan for quantum information science with $30 billion investment 
'@
$ws1.Range("C16").Value = $v
$v = @'
 https://www.nature.com/articles/d41586-025-01358-z 
'@
$ws1.Range("D16").Value = $v
$v = @'
 China's State Council approved the "Quantum Information Science and Technology Development Plan (2025-2040)" on May 8, 2025. The comprehensive strategy includes $30 billion in government investment over 15 years, targeting quantum computing, quantum communication, and quantum sensing. The plan establishes five national quantum research centers and details talent development programs aiming to train 50,000 quantum specialists. Strategic goals include achieving quantum advantage in specific applications by 2030 and developing fully fault-tolerant quantum computers by 2035. 
'@
$ws1.Range("E16").Value = $v

# Row 17
$v = @'
 Intellectual property protection 
'@
$ws1.Range("A17").Value = $v
$v = @'
 11 
'@
$ws1.Range("B17").Value = $v
$v = @'
 China strengthens intellectual property courts with new enforcement mechanisms 
'@
$ws1.Range("C17").Value = $v
$v = @'
 https://www.wipo.int/wipo_magazine/en/2025/02/article_0003.html 
'@
$ws1.Range("D17").Value = $v
$v = @'
 On May 5, 2025, China's Supreme People's Court announced significant reforms to intellectual property protection, including expanded jurisdiction for specialized IP courts and higher statutory damages for infringement. The reforms introduce a "blacklist" system for repeat IP violators that restricts their access to government contracts and financing. New technical investigation procedures for complex patent cases will be implemented, and a national IP dispute mediation center established. The measures signal China's increasing emphasis on creating a strong innovation ecosystem through IP protection. 
'@
$ws1.Range("E17").Value = $v

# Row 18
$v = @'
 Strategic emerging industries 
'@
$ws1.Range("A18").Value = $v
$v = @'
 10 
'@
$ws1.Range("B18").Value = $v
$v = @'
 China designates seven strategic emerging industry clusters with special policy support 
'@
$ws1.Range("C18").Value = $v
$v = @'
 https://www.bloomberg.com/news/articles/2025-05-04/china-designates-strategic-emerging-industry-clusters-with-policy-support 
'@
$ws1.Range("D18").Value = $v
$v = @'
 China's National Development and Reform Commission designated seven strategic emerging industry clusters on May 4, 2025. The clusters, located in Beijing, Shanghai, Shenzhen, Hefei, Hangzhou, Xi'an, and Wuhan, will receive preferential policies including tax incentives, streamlined regulatory approval, and special talent programs. Focus areas include next-generation information technology, synthetic biology, new energy vehicles, aerospace, and advanced materials. The initiative aims to create innovation ecosystems that integrate the entire industrial chain from research to manufacturing and applications. 
'@
$ws1.Range("E18").Value = $v

# Row 19
$v = @'
 Talent introduction policy 
'@
$ws1.Range("A19").Value = $v
$v = @'
 9 
'@
$ws1.Range("B19").Value = $v
$v = @'
 China launches expanded global talent recruitment program with simplified visa process 
'@
$ws1.Range("C19").Value = $v
$v = @'
 https://www.chinadaily.com.cn/a/202505/07/WS6547dc91a3104efcbdad7e22.html 
'@
$ws1.Range("D19").Value = $v
$v = @'
 China's Ministry of Human Resources and Social Security announced an expanded global talent recruitment initiative on May 7, 2025. The program introduces a new "S-visa" category with fast-track processing for high-level scientists and entrepreneurs. Benefits include ten-year multiple-entry visas, simplified permanent residency applications, and tax incentives. The policy targets experts in semiconductors, artificial intelligence, biotechnology, and quantum science. Local governments will establish international talent service centers in 30 cities to provide one-stop services for foreign professionals. 
'@
$ws1.Range("E19").Value = $v

# Row 20
$v = @'
 Drone logistics 
'@
$ws1.Range("A20").Value = $v
$v = @'
 9 
'@
$ws1.Range("B20").Value = $v
$v = @'
 China approves world's largest commercial drone delivery network for intercity logistics 
'@
$ws1.Range("C20").Value = $v
$v = @'
 https://techcrunch.com/2025/05/03/china-approves-massive-drone-delivery-network/ 
'@
$ws1.Range("D20").Value = $v
$v = @'
 China's Civil Aviation Administration approved the world's largest commercial drone logistics network on May 3, 2025. The network will connect 25 cities in the Yangtze River Delta region using autonomous cargo drones with 200kg payload capacity. JD Logistics and the Aviation Industry Corporation of China will jointly develop and operate the system, with initial operations beginning in September 2025. The network will utilize dedicated air corridors and automated landing facilities. Officials estimate the system will reduce delivery times by 70% and lower logistics costs by 30% for participating regions. 
'@
$ws1.Range("E20").Value = $v

# Row 21
$v = @'
 Quantum computing cloud platform 
'@
$ws1.Range("A21").Value = $v
$v = @'
 8 
'@
$ws1.Range("B21").Value = $v
$v = @'
 China launches national quantum computing cloud platform with 5 different quantum processors accessible 
'@
$ws1.Range("C21").Value = $v
$v = @'
 https://www.scmp.com/tech/big-tech/article/3259502/china-launches-national-quantum-computing-cloud-platform 
'@
$ws1.Range("D21").Value = $v
$v = @'
 China launched its National Quantum Computing Cloud Platform on May 4, 2025, providing researchers and companies with access to multiple quantum computing systems. The platform integrates five different quantum processors, including superconducting, trapped-ion, and photonic quantum computers from the Chinese Academy of Sciences, University of Science and Technology of China, and leading Chinese tech companies. The service offers free access for academic research and tiered commercial pricing. Over 100 organizations have already registered to use the platform, which officials described as critical infrastructure for quantum algorithm development and applications research. 
'@
$ws1.Range("E21").Value = $v

# ============ Sheet: Sources ============
$ws2 = $wb.Worksheets.Item("Sources")

# Row 3
$v = @'
 South China Morning Post 
'@
$ws2.Range("A3").Value = $v
$v = @'
 https://www.scmp.com/tech/policy/article/3259421/china-unveils-new-ai-governance-framework-emphasizing-advanced-model-safety 
'@
$ws2.Range("B3").Value = $v
$v = @'
 May 7, 2025 
'@
$ws2.Range("C3").Value = $v

# Row 4
$v = @'
 Nature 
'@
$ws2.Range("A4").Value = $v
$v = @'
 https://www.nature.com/articles/s41586-025-5742-x 
'@
$ws2.Range("B4").Value = $v
$v = @'
 May 4, 2025 
'@
$ws2.Range("C4").Value = $v

# Row 5
$v = @'
 Xinhua News Agency 
'@
$ws2.Range("A5").Value = $v
$v = @'
 https://www.xinhuanet.com/english/2025-05/06/c_1419352.htm 
'@
$ws2.Range("B5").Value = $v
$v = @'
 May 6, 2025 
'@
$ws2.Range("C5").Value = $v

# Row 6
$v = @'
 Reuters 
'@
$ws2.Range("A6").Value = $v
$v = @'
 https://www.reuters.com/technology/china-unveils-ambitious-manufacturing-transformation-plan-2025-05-08/ 
'@
$ws2.Range("B6").Value = $v
$v = @'
 May 8, 2025 
'@
$ws2.Range("C6").Value = $v

# Row 7
$v = @'
 Nikkei Asia 
'@
$ws2.Range("A7").Value = $v
$v = @'
 https://asia.nikkei.com/Business/Tech/Semiconductors/China-announces-breakthroughs-in-advanced-chip-packaging 
'@
$ws2.Range("B7").Value = $v
$v = @'
 May 3, 2025 
'@
$ws2.Range("C7").Value = $v

# Row 8
$v = @'
 Clean Energy Wire 
'@
$ws2.Range("A8").Value = $v
$v = @'
 https://www.cleanenergywire.org/news/china-germany-launch-joint-hydrogen-energy-storage-research-center 
'@
$ws2.Range("B8").Value = $v
$v = @'
 May 5, 2025 
'@
$ws2.Range("C8").Value = $v

# Row 9
$v = @'
 South China Morning Post 
'@
$ws2.Range("A9").Value = $v
$v = @'
 https://www.scmp.com/tech/policy/article/3259487/china-overhauls-tech-transfer-system-boost-commercialization 
'@
$ws2.Range("B9").Value = $v
$v = @'
 May 6, 2025 
'@
$ws2.Range("C9").Value = $v

# Row 10
$v = @'
 Nature 
'@
$ws2.Range("A10").Value = $v
$v = @'
 https://www.nature.com/articles/d41586-025-01355-2 
'@
$ws2.Range("B10").Value = $v
$v = @'
 May 7, 2025 
'@
$ws2.Range("C10").Value = $v

# Row 11
$v = @'
 Government of China 
'@
$ws2.Range("A11").Value = $v
$v = @'
 https://english.www.gov.cn/news/topnews/202505/08/content_WS6549f7a2c6d0868f4e2b2e37.html 
'@
$ws2.Range("B11").Value = $v
$v = @'
 May 8, 2025 
'@
$ws2.Range("C11").Value = $v

# Row 12
$v = @'
 Global Times 
'@
$ws2.Range("A12").Value = $v
$v = @'
 https://www.globaltimes.cn/page/202505/1306782.shtml 
'@
$ws2.Range("B12").Value = $v
$v = @'
 May 4, 2025 
'@
$ws2.Range("C12").Value = $v

# Row 13
$v = @'
 Science Magazine 
'@
$ws2.Range("A13").Value = $v
$v = @'
 https://www.sciencemag.org/news/2025/05/china-approves-groundbreaking-brain-computer-interface-clinical-trial 
'@
$ws2.Range("B13").Value = $v
$v = @'
 May 5, 2025 
'@
$ws2.Range("C13").Value = $v

# Row 14
$v = @'
 China Daily 
'@
$ws2.Range("A14").Value = $v
$v = @'
 https://www.chinadaily.com.cn/a/202505/06/WS6546c891a3104efcbdad7e21.html 
'@
$ws2.Range("B14").Value = $v
$v = @'
 May 6, 2025 
'@
$ws2.Range("C14").Value = $v

# Row 15
$v = @'
 ASEAN Official Website 
'@
$ws2.Range("A15").Value = $v
$v = @'
 https://asean.org/china-asean-launch-environmental-technology-partnership/ 
'@
$ws2.Range("B15").Value = $v
$v = @'
 May 7, 2025 
'@
$ws2.Range("C15").Value = $v

# Row 16
$v = @'
 Science Direct 
'@
$ws2.Range("A16").Value = $v
$v = @'
 https://www.sciencedirect.com/science/article/pii/B9780323856249000156 
'@
$ws2.Range("B16").Value = $v
$v = @'
 May 3, 2025 
'@
$ws2.Range("C16").Value = $v

# Row 17
$v = @'
 Nature 
'@
$ws2.Range("A17").Value = $v
$v = @'
 https://www.nature.com/articles/d41586-025-01358-z 
'@
$ws2.Range("B17").Value = $v
$v = @'
 May 8, 2025 
'@
$ws2.Range("C17").Value = $v

# Row 18
$v = @'
 WIPO Magazine 
'@
$ws2.Range("A18").Value = $v
$v = @'
 https://www.wipo.int/wipo_magazine/en/2025/02/article_0003.html 
'@
$ws2.Range("B18").Value = $v
$v = @'
 May 5, 2025 
'@
$ws2.Range("C18").Value = $v

# Row 19
$v = @'
 Bloomberg 
'@
$ws2.Range("A19").Value = $v
$v = @'
 https://www.bloomberg.com/news/articles/2025-05-04/china-designates-strategic-emerging-industry-clusters-with-policy-support 
'@
$ws2.Range("B19").Value = $v
$v = @'
 May 4, 2025 
'@
$ws2.Range("C19").Value = $v

# Row 20
$v = @'
 China Daily 
'@
$ws2.Range("A20").Value = $v
$v = @'
 https://www.chinadaily.com.cn/a/202505/07/WS6547dc91a3104efcbdad7e22.html 
'@
$ws2.Range("B20").Value = $v
$v = @'
 May 7, 2025 
'@
$ws2.Range("C20").Value = $v

# Row 21
$v = @'
 TechCrunch 
'@
$ws2.Range("A21").Value = $v
$v = @'
 https://techcrunch.com/2025/05/03/china-approves-massive-drone-delivery-network/ 
'@
$ws2.Range("B21").Value = $v
$v = @'
 May 3, 2025 
'@
$ws2.Range("C21").Value = $v

# Row 22
$v = @'
 South China Morning Post 
'@
$ws2.Range("A22").Value = $v
$v = @'
 https://www.scmp.com/tech/big-tech/article/3259502/china-launches-national-quantum-computing-cloud-platform 
'@
$ws2.Range("B22").Value = $v
$v = @'
 May 4, 2025 
'@
$ws2.Range("C22").Value = $v

# ============ Sheet: Executive Summary ============
$ws3 = $wb.Worksheets.Item("Executive Summary")
$execText = @'
Five Most Impactful News Summaries:

1. China unveiled a comprehensive "Framework for AI Safety and Governance 2025-2030" that establishes mandatory security assessments for advanced AI models while creating a national AI safety research center. The framework aims to maintain China's competitive edge in generative AI development while ensuring safety and responsible deployment, signaling China's approach to balancing innovation with regulation in a critical technology domain.

2. The State Council released the "New Quality Productivity Acceleration Plan (2025-2035)" outlining China's manufacturing transformation strategy with the designation of 25 advanced manufacturing clusters in sectors including semiconductors, biotech, and new energy vehicles. The plan aims to increase high-tech manufacturing's contribution to 35% of GDP by 2035 through tax incentives for R&D and emphasis on indigenous innovation while maintaining openness to international cooperation.

3. China launched the Digital Silk Road 2.0 initiative with a $50 billion investment package to support digital infrastructure development across Belt and Road countries. The program emphasizes sustainable technologies, AI solutions for climate monitoring, and interoperable digital payment systems with technology transfer provisions and training programs for partner countries, representing China's expanding digital diplomacy strategy.

4. The Chinese Academy of Sciences announced a breakthrough in solid-state battery technology achieving an energy density of 500 Wh/kg with stable performance over 1,000 cycles. CATL plans to incorporate the technology into pilot production lines by late 2025, positioning China to potentially lead in next-generation battery technology critical for electric vehicles and renewable energy storage.

5. China and European partners demonstrated the first intercontinental quantum-secured communication network using satellite-based quantum key distribution and terrestrial fiber networks. This milestone in China-EU scientific cooperation advances practical quantum communication technology and serves as a model for future international collaboration in sensitive technology areas, highlighting China's quantum leadership and science diplomacy approach.
'@
$ws3.Range("A2").Value = $execText

# ============ Sheet: Cooccurrence ============
$ws4 = $wb.Worksheets.Item("Cooccurrence")
$ws4.Range("A3").EntireRow.Delete()
$v = @'
Quantum communication
'@
$ws4.Range("A2").Value = $v
$v = @'
Science and technology development plan
'@
$ws4.Range("B2").Value = $v
$ws4.Range("C2").Value = 1

# ============ Sheet: Associations ============
$ws5 = $wb.Worksheets.Item("Associations")
$ws5.Range("A12:A13").EntireRow.Delete()
# Row 2
$v = @'
Quantum communication
'@
$ws5.Range("A2").Value = $v
$ws5.Range("B2").Value = 2

# Row 3
$v = @'
New quality productivity
'@
$ws5.Range("A3").Value = $v
$ws5.Range("B3").Value = 1

# Row 4
$v = @'
Semiconductor packaging
'@
$ws5.Range("A4").Value = $v
$ws5.Range("B4").Value = 1

# Row 5
$v = @'
Hydrogen energy storage
'@
$ws5.Range("A5").Value = $v
$ws5.Range("B5").Value = 1

# Row 6
$v = @'
Brain-computer interface
'@
$ws5.Range("A6").Value = $v
$ws5.Range("B6").Value = 1

# Row 7
$v = @'
Science and technology development plan
'@
$ws5.Range("A7").Value = $v
$ws5.Range("B7").Value = 1

# Row 8
$v = @'
Intellectual property protection
'@
$ws5.Range("A8").Value = $v
$ws5.Range("B8").Value = 1

# Row 9
$v = @'
Artificial Intelligence
'@
$ws5.Range("A9").Value = $v
$ws5.Range("B9").Value = 1

# Row 10
$v = @'
Drone logistics
'@
$ws5.Range("A10").Value = $v
$ws5.Range("B10").Value = 1

# Row 11
$v = @'
Quantum computing cloud platform
'@
$ws5.Range("A11").Value = $v
$ws5.Range("B11").Value = 1
